$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (73) below the existing last row (72), mirroring its
# column layout / style (date formatting on column D).
$row = 73

$ws.Cells.Item($row, 1).Value  = 4
$ws.Cells.Item($row, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value  = "Los Lagos"

$ws.Cells.Item($row, 4).Value  = 44939
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat

$ws.Cells.Item($row, 5).Value  = 10
$ws.Cells.Item($row, 6).Value  = "Fruta"
$ws.Cells.Item($row, 7).Value  = 100103
$ws.Cells.Item($row, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item($row, 9).Value  = 100103003
$ws.Cells.Item($row, 10).Value = "Damasco"
$ws.Cells.Item($row, 11).Value = "Modesto"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 600
$ws.Cells.Item($row, 14).Value = 20000
$ws.Cells.Item($row, 15).Value = 21000
$ws.Cells.Item($row, 16).Value = 20500
$ws.Cells.Item($row, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item($row, 19).Value = 1281
$ws.Cells.Item($row, 20).Value = 16
